$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp text (A1)
$ws.Range("A1").Value = "Datos actualizados a 31 de Julio de 2020 a las 21:34"

# Row 4 - Estados Unidos: refresh totals
$ws.Range("B4").Value = 4675916
$ws.Range("C4").Value = 40931
$ws.Range("D4").Value = 2299080
$ws.Range("E4").Value = 2220709
$ws.Range("G4").Value = 842
$ws.Range("H4").Value = 156127

# Row 6 - India: refresh totals
$ws.Range("B6").Value = 1696780
$ws.Range("C6").Value = 57430
$ws.Range("D6").Value = 1095647
$ws.Range("E6").Value = 564582
$ws.Range("G6").Value = 765
$ws.Range("H6").Value = 36551

# Row 21 - Alemania: refresh totals
$ws.Range("B21").Value = 210333
$ws.Range("C21").Value = 680
$ws.Range("E21").Value = 8810
$ws.Range("G21").Value = 2
$ws.Range("H21").Value = 9223

# Row 65 - Uzbekistan: refresh totals
$ws.Range("B65").Value = 24009
$ws.Range("C65").Value = 738
$ws.Range("D65").Value = 14464
$ws.Range("E65").Value = 9404
$ws.Range("G65").Value = 5
$ws.Range("H65").Value = 141

# Namibia moves up above Lituania and Estonia with a fresh data snapshot;
# Lituania and Estonia each shift down one row, keeping their prior values.
$ws.Range("A126").Value = "Namibia"
$ws.Range("B126").Value = 2129
$ws.Range("C126").Value = 77
$ws.Range("D126").Value = 166
$ws.Range("E126").Value = 1953
$ws.Range("H126").Value = 10

$ws.Range("A127").Value = "Lituania"
$ws.Range("B127").Value = 2075
$ws.Range("C127").Value = 13
$ws.Range("D127").Value = 1644
$ws.Range("E127").Value = 351
$ws.Range("H127").Value = 80

$ws.Range("A128").Value = "Estonia"
$ws.Range("B128").Value = 2064
$ws.Range("C128").Value = 13
$ws.Range("D128").Value = 1930
$ws.Range("E128").Value = 65
$ws.Range("H128").Value = 69
